$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume number + date range ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Murder row (row 14): counts suppressed to text placeholders "0" / "***.*" ---
# Force text storage (Excel would otherwise auto-coerce "0" back to a number),
# then copy the General-format/right-aligned style already used by the other
# placeholder cells (e.g. N22) so the format matches exactly.
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("C14:E14").PasteSpecial(-4122)

# --- Remaining updated statistics cells (rows 14-30) ---
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 0
$ws.Range("L14").Value = -9.756097560975
$ws.Range("N14").Value = -83.027522935779
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -33.333333333333
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = -16.666666666666
$ws.Range("I15").Value = 83
$ws.Range("J15").Value = 122
$ws.Range("K15").Value = -31.967213114754
$ws.Range("L15").Value = -7.777777777777
$ws.Range("M15").Value = -37.593984962406
$ws.Range("N15").Value = -74.143302180685
$ws.Range("C16").Value = 40
$ws.Range("D16").Value = 51
$ws.Range("E16").Value = -21.56862745098
$ws.Range("F16").Value = 163
$ws.Range("H16").Value = -19.306930693069
$ws.Range("I16").Value = 1201
$ws.Range("J16").Value = 1344
$ws.Range("K16").Value = -10.639880952381
$ws.Range("L16").Value = 11.410018552875
$ws.Range("M16").Value = -18.410326086956
$ws.Range("N16").Value = -80.197856553998
$ws.Range("C17").Value = 57
$ws.Range("D17").Value = 63
$ws.Range("E17").Value = -9.523809523809
$ws.Range("F17").Value = 243
$ws.Range("G17").Value = 274
$ws.Range("H17").Value = -11.313868613138
$ws.Range("I17").Value = 1914
$ws.Range("J17").Value = 2001
$ws.Range("K17").Value = -4.347826086956
$ws.Range("L17").Value = 9.371428571428
$ws.Range("M17").Value = 52.631578947368
$ws.Range("N17").Value = -49.750590706222
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = -40.54054054054
$ws.Range("F18").Value = 104
$ws.Range("G18").Value = 129
$ws.Range("H18").Value = -19.37984496124
$ws.Range("I18").Value = 963
$ws.Range("J18").Value = 1150
$ws.Range("K18").Value = -16.260869565217
$ws.Range("L18").Value = 12.5
$ws.Range("M18").Value = 12.763466042154
$ws.Range("N18").Value = -86.23499142367
$ws.Range("C19").Value = 156
$ws.Range("D19").Value = 139
$ws.Range("E19").Value = 12.230215827338
$ws.Range("F19").Value = 621
$ws.Range("G19").Value = 536
$ws.Range("H19").Value = 15.858208955223
$ws.Range("I19").Value = 4028
$ws.Range("J19").Value = 4106
$ws.Range("K19").Value = -1.899659035557
$ws.Range("L19").Value = 30.736773774748
$ws.Range("M19").Value = 37.803626411221
$ws.Range("N19").Value = -42.159678345778
$ws.Range("C20").Value = 26
$ws.Range("E20").Value = 36.842105263157
$ws.Range("F20").Value = 99
$ws.Range("G20").Value = 105
$ws.Range("H20").Value = -5.714285714285
$ws.Range("I20").Value = 868
$ws.Range("J20").Value = 827
$ws.Range("K20").Value = 4.957678355501
$ws.Range("L20").Value = 43.946932006633
$ws.Range("M20").Value = 138.461538461538
$ws.Range("N20").Value = -85.399495374264
$ws.Range("C21").Value = 303
$ws.Range("D21").Value = 312
$ws.Range("E21").Value = -2.884615384615
$ws.Range("F21").Value = 1243
$ws.Range("G21").Value = 1261
$ws.Range("H21").Value = -1.42743854084
$ws.Range("I21").Value = 9094
$ws.Range("J21").Value = 9582
$ws.Range("K21").Value = -5.092882487998
$ws.Range("L21").Value = 21.269502600346
$ws.Range("M21").Value = 29.304706384188
$ws.Range("N21").Value = -70.004617718846
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = 9.090909090909
$ws.Range("I22").Value = 183
$ws.Range("J22").Value = 199
$ws.Range("K22").Value = -8.040201005025
$ws.Range("L22").Value = 29.787234042553
$ws.Range("M22").Value = 22.818791946308
$ws.Range("C23").Value = 20
$ws.Range("D23").Value = 26
$ws.Range("E23").Value = -23.076923076923
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 109
$ws.Range("H23").Value = 1.834862385321
$ws.Range("I23").Value = 800
$ws.Range("J23").Value = 827
$ws.Range("K23").Value = -3.264812575574
$ws.Range("L23").Value = 3.626943005181
$ws.Range("M23").Value = 52.091254752851
$ws.Range("C24").Value = 242
$ws.Range("D24").Value = 358
$ws.Range("E24").Value = -32.402234636871
$ws.Range("F24").Value = 1100
$ws.Range("G24").Value = 1257
$ws.Range("H24").Value = -12.490055688146
$ws.Range("I24").Value = 8975
$ws.Range("J24").Value = 10052
$ws.Range("K24").Value = -10.714285714285
$ws.Range("L24").Value = 17.828541420506
$ws.Range("M24").Value = 54.262633207287
$ws.Range("C25").Value = 98
$ws.Range("D25").Value = 66
$ws.Range("E25").Value = 48.484848484848
$ws.Range("F25").Value = 352
$ws.Range("G25").Value = 313
$ws.Range("H25").Value = 12.460063897763
$ws.Range("I25").Value = 2920
$ws.Range("J25").Value = 2897
$ws.Range("K25").Value = 0.793924749741
$ws.Range("L25").Value = 12.91569992266
$ws.Range("M25").Value = -16.499857020303
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = -5.882352941176
$ws.Range("I26").Value = 156
$ws.Range("J26").Value = 190
$ws.Range("K26").Value = -17.894736842105
$ws.Range("L26").Value = -2.5
$ws.Range("C27").Value = 6
$ws.Range("E27").Value = -45.454545454545
$ws.Range("F27").Value = 40
$ws.Range("H27").Value = -13.043478260869
$ws.Range("I27").Value = 384
$ws.Range("J27").Value = 436
$ws.Range("K27").Value = -11.926605504587
$ws.Range("L27").Value = -6.569343065693
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 18
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 106
$ws.Range("J28").Value = 133
$ws.Range("K28").Value = -20.300751879699
$ws.Range("L28").Value = -33.75
$ws.Range("M28").Value = -19.083969465648
$ws.Range("N28").Value = -79.576107899807
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -25
$ws.Range("F29").Value = 18
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = 20
$ws.Range("I29").Value = 95
$ws.Range("J29").Value = 111
$ws.Range("K29").Value = -14.414414414414
$ws.Range("L29").Value = -33.098591549295
$ws.Range("M29").Value = -17.391304347826
$ws.Range("N29").Value = -80.083857442348
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 5
$ws.Range("H30").Value = -28.571428571428
$ws.Range("I30").Value = 47
$ws.Range("J30").Value = 69
$ws.Range("K30").Value = -31.884057971014
$ws.Range("L30").Value = -9.615384615384
